# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.104.12'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.573.60'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.577.24'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('E12').Value = '  +11.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.345'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.022.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.124.32'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.43%  '
$ws.Range('E17').Value = '  +3.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.574.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.457'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.58%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.161'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.989'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.01'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  +1.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.871'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '292.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.57%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '127.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.54%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.592'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.947.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.26%  '
